$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the leftover "_GoBack" bookmark from the abstract
#    paragraph (present in the original, gone in the revised file).
# ------------------------------------------------------------------
try {
    $gb = $d.Bookmarks("_GoBack")
    $gb.Delete()
} catch {
    # bookmark may already be gone - nothing to do
}

# ------------------------------------------------------------------
# 2. Re-flow the abstract paragraph's runs.  The visible text is
#    unchanged, but the revision consolidates the many small runs
#    that resulted from live editing/track-changes back into a
#    handful of runs.  Doing a Find/Replace with identical text over
#    each stretch forces Word to rebuild the run for that stretch,
#    which also removes the now-orphan bookmark markers that sit
#    inside the replaced range.  The "crops, and" sentence keeps its
#    own run (it is wrapped in proofErr markers flagging a grammar
#    check) so that boundary is left untouched.
# ------------------------------------------------------------------
$chunk1 = "Agroecologists are often concerned with the indirect effects of experiment treatments, management practices, and environmental gradients. However, standard statistical approaches using multiple regression or generalized linear models are not suited to quantify direct vs. indirect effects among a network of interactions. Path analysis, a type of structural equation modeling, has been increasingly appreciated in community ecology as an important tool for quantifying and evaluating indirect effects. Relatively user-friendly R packages have been developed that require only basic knowledge of R commands and linear models to implement, but they have not yet been extensively adopted in agricultural research. In this paper, we provide several examples from real agroecological experiments conducted in the Pacific Northwest "
$r1 = $d.Content
$r1.Find.Execute($chunk1, $true, $false, $false, $false, $false, $true, 1, $false, $chunk1, 2) | Out-Null

$chunk2 = " demonstrate the usefulness of path analysis and how to implement it. Importantly, we demonstrate that important biological inferences would be otherwise obfuscated had path analysis not been utilized as a statistical tool.  We urge other researches to attempt this approach and highlight several examples, such as quantification of biological control and crop yield, where path analysis should be considered a default tool for the evaluation of experimental outcomes."
$r2 = $d.Content
$r2.Find.Execute($chunk2, $true, $false, $false, $false, $false, $true, 1, $false, $chunk2, 2) | Out-Null

# ------------------------------------------------------------------
# 3. Add two new paragraphs right after the abstract paragraph:
#      - an empty paragraph
#      - a paragraph containing "Add text for example."
#    Both inherit the abstract paragraph's formatting
#    (firstLine indent 720, contextualSpacing 0).
# ------------------------------------------------------------------
$abstractPara = $d.Paragraphs.Item(7)
$insertPoint = $abstractPara.Range
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()
$insertPoint.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Item(9)
$newPara2.Range.InsertAfter("Add text for example.")

Write-Host "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
